$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set cells as Text format first to preserve values exactly (avoid numeric/date auto-conversion)
$cells = @(
    @{Addr='D2'; Val='65.386.67'}
    @{Addr='E2'; Val='  -0.77%  '}
    @{Addr='D3'; Val='2.936.22'}
    @{Addr='E3'; Val='  -2.58%  '}
    @{Addr='E4'; Val='  -0.04%  '}
    @{Addr='D5'; Val='568.30'}
    @{Addr='E5'; Val='  -2.88%  '}
    @{Addr='D6'; Val='158.72'}
    @{Addr='E6'; Val='  +1.86%  '}
    @{Addr='E7'; Val='  +0.08%  '}
    @{Addr='E8'; Val='  -0.38%  '}
    @{Addr='D9'; Val='2.931.49'}
    @{Addr='E9'; Val='  -2.64%  '}
    @{Addr='D10'; Val='6.69'}
    @{Addr='E10'; Val='  -3.38%  '}
    @{Addr='E11'; Val='  -3.67%  '}
    @{Addr='E12'; Val='  +1.60%  '}
    @{Addr='E13'; Val='  -2.29%  '}
    @{Addr='D14'; Val='34.27'}
    @{Addr='E14'; Val='  -0.94%  '}
    @{Addr='E15'; Val='  -0.70%  '}
    @{Addr='D16'; Val='65.410.87'}
    @{Addr='E16'; Val='  -0.75%  '}
    @{Addr='D17'; Val='3.425.76'}
    @{Addr='E17'; Val='  -2.51%  '}
    @{Addr='D18'; Val='6.97'}
    @{Addr='E18'; Val='  +0.12%  '}
    @{Addr='D19'; Val='2.937.34'}
    @{Addr='E19'; Val='  -2.34%  '}
    @{Addr='D20'; Val='15.68'}
    @{Addr='E20'; Val='  +13.40%  '}
    @{Addr='D21'; Val='444.79'}
    @{Addr='E21'; Val='  -4.10%  '}
    @{Addr='E22'; Val='  +0.73%  '}
    @{Addr='D23'; Val='7.25'}
    @{Addr='E23'; Val='  -1.73%  '}
    @{Addr='D24'; Val='82.16'}
    @{Addr='E24'; Val='  +0.22%  '}
    @{Addr='E25'; Val='  -1.40%  '}
    @{Addr='E26'; Val='  -3.30%  '}
    @{Addr='D27'; Val='10.06'}
    @{Addr='E27'; Val='  -5.88%  '}
    @{Addr='E28'; Val='  +0.10%  '}
    @{Addr='D29'; Val='8.06'}
    @{Addr='E29'; Val='  +1.05%  '}
    @{Addr='D30'; Val='2.36'}
    @{Addr='E30'; Val='  -1.07%  '}
    @{Addr='D31'; Val='2.58'}
    @{Addr='E31'; Val='  -1.48%  '}
    @{Addr='E32'; Val='  -4.31%  '}
    @{Addr='D33'; Val='27.09'}
    @{Addr='E33'; Val='  +0.29%  '}
    @{Addr='E34'; Val='  -0.24%  '}
    @{Addr='E35'; Val='  +0.03%  '}
    @{Addr='D36'; Val='0.972'}
    @{Addr='E36'; Val='  -2.69%  '}
    @{Addr='D37'; Val='5.72'}
    @{Addr='E37'; Val='  -1.64%  '}
    @{Addr='D38'; Val='49.70'}
    @{Addr='E38'; Val='  +0.93%  '}
    @{Addr='D39'; Val='45.18'}
    @{Addr='E39'; Val='  +1.88%  '}
    @{Addr='E40'; Val='  -9.25%  '}
    @{Addr='D41'; Val='0.301'}
    @{Addr='E41'; Val='  -0.72%  '}
    @{Addr='E42'; Val='  -1.99%  '}
    @{Addr='D43'; Val='2.83'}
    @{Addr='E43'; Val='  -6.65%  '}
    @{Addr='D44'; Val='8.46'}
    @{Addr='E44'; Val='  +0.04%  '}
    @{Addr='D45'; Val='383.51'}
    @{Addr='E45'; Val='  -3.46%  '}
    @{Addr='D46'; Val='0.0351'}
    @{Addr='E46'; Val='  -0.79%  '}
    @{Addr='D47'; Val='2.700.52'}
    @{Addr='E47'; Val='  -3.46%  '}
    @{Addr='D48'; Val='133.52'}
    @{Addr='E48'; Val='  -0.31%  '}
    @{Addr='E49'; Val='  +0.01%  '}
    @{Addr='E50'; Val='  +4.26%  '}
    @{Addr='D51'; Val='23.43'}
    @{Addr='E51'; Val='  -0.47%  '}
)

foreach ($item in $cells) {
    $rng = $ws.Range($item.Addr)
    $rng.NumberFormat = '@'
    $rng.Value = $item.Val
}
